$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.813.92'
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").Value = '1.942.11'
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("E4").Value = '  -0.11%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '242.59'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -2.07%  '
$ws.Range("E6").Value = '  -0.06%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4889'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.09%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2959'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.39%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06890'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +1.04%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.44'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.58%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '106.28'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("D12").Value = '1.941.51'
$ws.Range("E12").Value = '  -0.77%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.07730'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("E14").Value = '  -1.20%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.6984'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.96%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '277.43'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -1.51%  '
$ws.Range("D17").Value = '30.809.39'
$ws.Range("E17").Value = '  -0.99%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000007695'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("E19").Value = '  -0.57%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").Value = '2.191.39'
$ws.Range("E21").Value = '  +0.17%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.471'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.32%  '
$ws.Range("E23").Value = '  -0.08%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '6.512'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.19%  '
$ws.Range("E25").Value = '  -2.27%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '167.96'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.53%  '
$ws.Range("E27").Value = '  -1.34%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.159'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.18%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.1046'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.90%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.391'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -3.55%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.553'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.78%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.561'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -4.27%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.367'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -3.36%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.04847'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -2.95%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.7521'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.91%  '
$ws.Range("E36").Value = '  -0.42%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.0000'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.01%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.730'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.22%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.01995'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -2.63%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.657'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.85%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '78.14'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +6.24%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '6.494'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +1.35%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.099'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.56%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.9056'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +2.41%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '108.01'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.22%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.4406'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.75%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.9988'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.23%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.749'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +3.72%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '991.78'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.68%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.1244'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.61%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '9.317'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.25%  '
